$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45188 -> 2023-09-19) for
# every data row (rows 2-110). The update bumps that date forward by one day
# (45188 -> 45189 => 2023-09-20) for all of them.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
